$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letters for the 20-column table (A..T)
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

# Full target data for rows 2-7 (row 1 is the header, untouched).
# Sending cluster / Ligand symbol / Receptor symbol / Target cluster use
# shared strings: "ECs", "FAPs", "sCs", "Anpep", "Sele".
$rows = @(
    @{ rownum = 2;  A="ECs";  B="Anpep"; C="Sele"; D="ECs";
       E=3; F=1; G=8.180121999999999; H=24.540366;
       I=0.06766391458781856; J=0.06766391458781855;
       K=3; L=1; M=7.879565666666667; N=23.638697;
       O=0.9977172793687663; P=0.9977172793687664;
       Q=64.45580846034466; R=580.102276143102;
       S=0.06750945677399892; T=0.06750945677399892 },

    @{ rownum = 3;  A="ECs";  B="Anpep"; C="Sele"; D="sCs";
       E=3; F=1; G=8.180121999999999; H=24.540366;
       I=0.06766391458781856; J=0.06766391458781855;
       K=1; L=0.3333333333333333; M=0.018028; N=0.054084;
       O=0.002282720631233623; P=0.002282720631233623;
       Q=0.147471239416; R=1.327241154744;
       S=0.0001544578138196432; T=0.0001544578138196431 },

    @{ rownum = 4;  A="FAPs"; B="Anpep"; C="Sele"; D="ECs";
       E=3; F=1; G=100.9737753333333; H=302.921326;
       I=0.8352297080366586; J=0.8352297080366585;
       K=3; L=1; M=7.879565666666667; N=23.638697;
       O=0.9977172793687663; P=0.9977172793687664;
       Q=795.629493350247; R=7160.665440152223;
       S=0.8333231119503041; T=0.8333231119503041 },

    @{ rownum = 5;  A="FAPs"; B="Anpep"; C="Sele"; D="sCs";
       E=3; F=1; G=100.9737753333333; H=302.921326;
       I=0.8352297080366586; J=0.8352297080366585;
       K=1; L=0.3333333333333333; M=0.018028; N=0.054084;
       O=0.002282720631233623; P=0.002282720631233623;
       Q=1.820355221709333; R=16.383196995384;
       S=0.001906596086354516; T=0.001906596086354516 },

    @{ rownum = 6;  A="sCs";  B="Anpep"; C="Sele"; D="ECs";
       E=3; F=1; G=11.739522; H=35.218566;
       I=0.09710637737552288; J=0.09710637737552287;
       K=3; L=1; M=7.879565666666667; N=23.638697;
       O=0.9977172793687663; P=0.9977172793687664;
       Q=92.50233449427802; R=832.5210104485021;
       S=0.09688471064446341; T=0.09688471064446341 },

    @{ rownum = 7;  A="sCs";  B="Anpep"; C="Sele"; D="sCs";
       E=3; F=1; G=11.739522; H=35.218566;
       I=0.09710637737552288; J=0.09710637737552287;
       K=1; L=0.3333333333333333; M=0.018028; N=0.054084;
       O=0.002282720631233623; P=0.002282720631233623;
       Q=0.211640102616; R=1.904760923544;
       S=0.000221666731059464; T=0.000221666731059464 }
)

foreach ($row in $rows) {
    $rn = $row.rownum
    foreach ($col in $cols) {
        $ws.Range("$col$rn").Value = $row[$col]
    }
}
